# Insert two new data rows (a new "Primera"/"Segunda" Acelga price pair)
# right before the current row 580, shifting the existing rows 580:662 down
# to 582:664 to match the target diff (dimension grows from R662 to R664).
# Excel's Rows.Insert() shifts everything below down a row and carries the
# row-above formatting (incl. the date style on column D) onto the freshly
# inserted blank rows, same as a manual "insert row" in the Excel UI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("580:581").Insert()

# Row 580 - new "Primera" quality entry.
$ws.Range("A580").Value = 8
$ws.Range("B580").Value = "Terminal La Palmera de La Serena"
$ws.Range("C580").Value = "Coquimbo"
$ws.Range("D580").Value2 = 45077
$ws.Range("E580").Value = 4
$ws.Range("F580").Value = 100112009
$ws.Range("G580").Value = "Acelga"
$ws.Range("H580").Value = "Sin especificar"
$ws.Range("I580").Value = "Primera"
$ws.Range("J580").Value = 2000
$ws.Range("K580").Value = 500
$ws.Range("L580").Value = 600
$ws.Range("M580").Value = 550
$ws.Range("N580").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O580").Value = "Provincia del Elquí"
$ws.Range("P580").Value = 275
$ws.Range("Q580").Value = 2
$ws.Range("R580").Value = "Hortaliza"

# Row 581 - new "Segunda" quality entry (same date as row 580).
$ws.Range("A581").Value = 8
$ws.Range("B581").Value = "Terminal La Palmera de La Serena"
$ws.Range("C581").Value = "Coquimbo"
$ws.Range("D581").Value2 = 45077
$ws.Range("E581").Value = 4
$ws.Range("F581").Value = 100112009
$ws.Range("G581").Value = "Acelga"
$ws.Range("H581").Value = "Sin especificar"
$ws.Range("I581").Value = "Segunda"
$ws.Range("J581").Value = 1360
$ws.Range("K581").Value = 400
$ws.Range("L581").Value = 450
$ws.Range("M581").Value = 425
$ws.Range("N581").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O581").Value = "Provincia del Elquí"
$ws.Range("P581").Value = 212
$ws.Range("Q581").Value = 2
$ws.Range("R581").Value = "Hortaliza"
